# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Vega Monumental Concepción - Piña"
# at row 195, pushing the existing rows 195:301 down to 196:302.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 195 (existing rows 195-301 shift down to 196-302).
$ws.Rows("195:195").Insert()

# Populate the new row 195 with the latest weekly record.
$ws.Cells.Item(195, 1).Value  = 11
$ws.Cells.Item(195, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(195, 3).Value  = "Bíobío"
$ws.Cells.Item(195, 4).Value  = 45176
$ws.Cells.Item(195, 5).Value  = 8
$ws.Cells.Item(195, 6).Value  = "Fruta"
$ws.Cells.Item(195, 7).Value  = 100108
$ws.Cells.Item(195, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(195, 9).Value  = 100108005
$ws.Cells.Item(195, 10).Value = "Piña"
$ws.Cells.Item(195, 11).Value = "Caramelo"
$ws.Cells.Item(195, 12).Value = "Segunda"
$ws.Cells.Item(195, 13).Value = 270
$ws.Cells.Item(195, 14).Value = 19000
$ws.Cells.Item(195, 15).Value = 20000
$ws.Cells.Item(195, 16).Value = 19556
$ws.Cells.Item(195, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(195, 18).Value = "Ecuador"
$ws.Cells.Item(195, 19).Value = 1397
$ws.Cells.Item(195, 20).Value = 14
